# Apply the authored edit to the АИ (AHP) workbook:
#  - Rename the only sheet "Лист2" -> "Лист1" and give it a fresh sheetId
#    (Excel mints a new internal sheetId when a sheet is duplicated; we
#    replicate that by copying the sheet, dropping the original, and
#    renaming the copy).
#  - Swap the 2nd and 4th concept names (Recipe Manager <-> FitLife) that
#    live in B13/C13 and B15/C15 - every other "Recipe Manager"/"FitLife"
#    label on the sheet is produced by formulas referencing these cells,
#    so they recompute automatically.
#  - Update one raw pairwise-comparison input (N30) from 9 to 3.
#  - Update the view: zoom to 125%, and move the active selection to D85.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- give the sheet a new identity (sheetId) like a genuine duplicate ---
$ws.Copy($null, $ws)
$excel.DisplayAlerts = $false
$ws.Delete()
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Лист1"
$wb.Worksheets.Item(1).Activate()

# --- swap the "Recipe Manager" / "FitLife" concept rows ---
# (.Value2 is used for reading - .Value round-trips through a
# PowerShell variable incorrectly in this host)
$b13 = $ws.Cells.Item(13, 2).Value2
$b15 = $ws.Cells.Item(15, 2).Value2
$c13 = $ws.Cells.Item(13, 3).Value2
$c15 = $ws.Cells.Item(15, 3).Value2

$ws.Cells.Item(13, 2).Value = $b15
$ws.Cells.Item(15, 2).Value = $b13
$ws.Cells.Item(13, 3).Value = $c15
$ws.Cells.Item(15, 3).Value = $c13

# --- correct the pairwise comparison input for "practical usefulness" ---
$ws.Cells.Item(30, 14).Value = 3

# --- view changes: zoom to 125%, select D85 ---
$excel.ActiveWindow.Zoom = 125
$ws.Range("D85").Select()
